$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This template gains two new specification columns:
#   - "Operating Weight Range (kg)"  -> inserted before "Bucket Capacity (m³)"
#     (becomes column G)
#   - "Rated Power ISO14396 (kW)"    -> inserted before "Implement Circuit (MPa)"
#     (becomes column R, after the G insertion has shifted everything right)
# Inserting whole columns automatically shifts the existing headers, row 2
# data and column widths to the right, so we only need to fill in the two
# brand-new columns afterwards.
# ---------------------------------------------------------------------------

# 1) Insert the "Operating Weight Range (kg)" column at G (shifts old G:AB -> H:AC)
$ws.Columns("G").Insert()

# 2) Insert the "Rated Power ISO14396 (kW)" column at R (shifts old R:AC -> S:AD)
$ws.Columns("R").Insert()

# 3) Populate the new "Operating Weight Range (kg)" column (G)
$ws.Range("G1").Value = "Operating Weight Range (kg)"
$ws.Range("G2").Value = 4000
$ws.Columns("G").ColumnWidth = 25

# 4) Populate the new "Rated Power ISO14396 (kW)" column (R)
$ws.Range("R1").Value = "Rated Power ISO14396 (kW)"
$ws.Range("R2").Value = 21.2
$ws.Columns("R").ColumnWidth = 25

Write-Host "Inserted Operating Weight Range (kg) and Rated Power ISO14396 (kW) columns"
